$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 60, shifting rows 60:148 down to 61:149
$ws.Rows.Item(60).Insert()

# Populate the new row 60 with the new data point
$ws.Range("A60").Value = 4
$ws.Range("B60").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C60").Value = "Los Lagos"
$ws.Range("D60").Value = 44894
$ws.Range("E60").Value = 10
$ws.Range("F60").Value = 100112052
$ws.Range("G60").Value = "Albahaca"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 60
$ws.Range("K60").Value = 8000
$ws.Range("L60").Value = 8000
$ws.Range("M60").Value = 8000
$ws.Range("N60").Value = "$/docena de matas"
$ws.Range("O60").Value = "Región Metropolitana"
$ws.Range("P60").Value = 1333
$ws.Range("Q60").Value = 6
$ws.Range("R60").Value = "Hortaliza"
